$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "= 3 * ( n * (n + 1) ) / 2   (Gaussche Summenformel)":
#    split the run that holds the trailing tab + "(Gaussche Summenformel)"
#    text so the label gets its own run with a smaller font size (9pt).
# ---------------------------------------------------------------------------
$gaussPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Gaussche Summenformel*") {
        $gaussPara = $p
        break
    }
}
if ($gaussPara -eq $null) { throw "Gaussche paragraph not found" }

$gaussXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="006A6676" w:rsidRDefault="006A6676" w:rsidP="00BC707C"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:tab/><w:t>= 3 * ( n * (n + 1) ) / 2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>(Gaussche Summenformel)</w:t></w:r></w:p>'
[void]$gaussPara.Range.InsertXML($gaussXml)

# ---------------------------------------------------------------------------
# 2) Paragraph "Die Drei in der ersten Zeile ...": change wording at the end
#    ("... bis zu einem" -> "... i-viele Iterationen macht.") and shrink the
#    whole paragraph to a 9pt font.
# ---------------------------------------------------------------------------
$dreiPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Die Drei in der ersten Zeile*") {
        $dreiPara = $p
        break
    }
}
if ($dreiPara -eq $null) { throw "'Die Drei' paragraph not found" }

$dreiXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00BC707C" w:rsidRDefault="006A6676" w:rsidP="00BC707C"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Die Drei in der ersten Zeile dieser Formeln kommt von den drei Operationen in der dritten Codezeile. Die aufsteigenden Zahlen der Summe sind eine Konsequenz der inner</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>en For-shleiffe die i-viele Iterationen macht.</w:t></w:r></w:p>'
$dreiPara.Range.InsertXML($dreiXml)

# ---------------------------------------------------------------------------
# 3) The (until now empty) paragraph right after it gains a single tab run.
# ---------------------------------------------------------------------------
$dreiParaIndex = $dreiPara.Index
$tabPara = $d.Paragraphs($dreiParaIndex + 1)

$tabXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="006A6676" w:rsidRDefault="006A6676" w:rsidP="00BC707C"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:tab/></w:r></w:p>'
$tabPara.Range.InsertXML($tabXml)
